# CT-3729 More dispatch letters for complaints
# Resize the letterhead table: narrower overall width, small negative
# left indent (table pulled slightly left of the margin), columns
# rebalanced, and explicit heights set on the address row and the
# trailing blank row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Table width: 10060 -> 9966 twips (values are expressed in points = twips/20)
$t.PreferredWidth = 9966 / 20.0

# New table indent from the left margin: -142 twips
$t.Rows.LeftIndent = -142 / 20.0

# Column widths (also propagates to each cell's tcW): 5751/4309 -> 5757/4209
$t.Columns.Item(1).Width = 5757 / 20.0
$t.Columns.Item(2).Width = 4209 / 20.0

# Explicit row heights on the address row (2) and the final blank row (3)
$t.Rows.Item(2).Height = 2022 / 20.0
$t.Rows.Item(3).Height = 235 / 20.0
